$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.0186
$ws.Range("A3").Value = 0.016
$ws.Range("A4").Value = 0.0121
$ws.Range("A5").Value = 0.0105
$ws.Range("A6").Value = 0.0087
$ws.Range("A7").Value = 0.0053

$ws.Range("A8:B9").ClearContents()

$ws.Range("D8").Select()
